$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SLAB")

# Insert two new columns before column D (shifts existing quarters from D:K to F:M)
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy number formats (date format / number format) from the adjacent column F into the new D:E columns
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Set values for the two new quarter columns (D, E)
$ws.Range("D7").Value = 43463
$ws.Range("E7").Value = 43372
$ws.Range("D8").Value = 215500
$ws.Range("E8").Value = 230200
$ws.Range("D9").Value = 85300
$ws.Range("E9").Value = 94600
$ws.Range("D10").Value = 130200
$ws.Range("E10").Value = 135600
$ws.Range("D12").Value = 62900
$ws.Range("E12").Value = 61100
$ws.Range("D17").Value = 197200
$ws.Range("E17").Value = 205100
$ws.Range("D18").Value = 18300
$ws.Range("E18").Value = 25100
$ws.Range("D20").Value = -200
$ws.Range("E20").Value = 2200
$ws.Range("D21").Value = 33000
$ws.Range("E21").Value = 46300
$ws.Range("D22").Value = 5000
$ws.Range("E22").Value = 4900
$ws.Range("D23").Value = 13100
$ws.Range("E23").Value = 22300
$ws.Range("D24").Value = -2100
$ws.Range("E24").Value = -300
$ws.Range("D26").Value = 15200
$ws.Range("E26").Value = 22600
$ws.Range("D27").Value = 15200
$ws.Range("E27").Value = 22600
$ws.Range("D29").Value = -100
$ws.Range("E29").Value = 5200
$ws.Range("D32").Value = 200
$ws.Range("E32").Value = -2200
$ws.Range("D33").Value = 15100
$ws.Range("E33").Value = 27800
$ws.Range("D35").Value = 15100
$ws.Range("E35").Value = 27800
$ws.Range("D38").Value = 43463
$ws.Range("E38").Value = 43372
$ws.Range("D41").Value = 197000
$ws.Range("E41").Value = 225300
$ws.Range("D42").Value = 416800
$ws.Range("E42").Value = 376600
$ws.Range("D43").Value = 73200
$ws.Range("E43").Value = 74600
$ws.Range("D44").Value = 75000
$ws.Range("E44").Value = 77600
$ws.Range("D45").Value = 64700
$ws.Range("E45").Value = 47000
$ws.Range("D46").Value = 826600
$ws.Range("E46").Value = 801100
$ws.Range("D48").Value = 139000
$ws.Range("E48").Value = 135600
$ws.Range("D49").Value = 568200
$ws.Range("E49").Value = 578300
$ws.Range("D52").Value = 90500
$ws.Range("E52").Value = 89400
$ws.Range("D54").Value = 1624400
$ws.Range("E54").Value = 1604300
$ws.Range("D57").Value = 41200
$ws.Range("E57").Value = 43600
$ws.Range("D59").Value = 103700
$ws.Range("E59").Value = 103700
$ws.Range("D60").Value = 144800
$ws.Range("E60").Value = 147300
$ws.Range("D61").Value = 354800
$ws.Range("E61").Value = 351500
$ws.Range("D62").Value = 57400
$ws.Range("E62").Value = 57800
$ws.Range("D66").Value = 557100
$ws.Range("E66").Value = 556600
$ws.Range("D72").Value = 961300
$ws.Range("E72").Value = 946200
$ws.Range("D76").Value = 1067300
$ws.Range("E76").Value = 1047700
$ws.Range("D80").Value = 43463
$ws.Range("E80").Value = 43372
$ws.Range("D81").Value = 15100
$ws.Range("E81").Value = 27800
$ws.Range("D83").Value = 14900
$ws.Range("E83").Value = 19100
$ws.Range("D89").Value = 28400
$ws.Range("E89").Value = 94300
$ws.Range("D91").Value = -6200
$ws.Range("E91").Value = -6900
$ws.Range("D94").Value = -47900
$ws.Range("E94").Value = -39100
$ws.Range("D100").Value = -8800
$ws.Range("E100").Value = -24800
$ws.Range("D102").Value = -28300
$ws.Range("E102").Value = 30500

# NA rows
$ws.Range("D47").Value = "NA"
$ws.Range("E47").Value = "NA"

# Zero rows
$ws.Range("D13:E13").Value = 0
$ws.Range("D14:E14").Value = 0
$ws.Range("D15:E15").Value = 0
$ws.Range("D25:E25").Value = 0
$ws.Range("D28:E28").Value = 0
$ws.Range("D30:E30").Value = 0
$ws.Range("D31:E31").Value = 0
$ws.Range("D34:E34").Value = 0
$ws.Range("D50:E50").Value = 0
$ws.Range("D51:E51").Value = 0
$ws.Range("D53:E53").Value = 0
$ws.Range("D58:E58").Value = 0
$ws.Range("D63:E63").Value = 0
$ws.Range("D64:E64").Value = 0
$ws.Range("D65:E65").Value = 0
$ws.Range("D68:E68").Value = 0
$ws.Range("D69:E69").Value = 0
$ws.Range("D70:E70").Value = 0
$ws.Range("D71:E71").Value = 0
$ws.Range("D73:E73").Value = 0
$ws.Range("D74:E74").Value = 0
$ws.Range("D75:E75").Value = 0
$ws.Range("D77:E77").Value = 0
$ws.Range("D84:E84").Value = 0
$ws.Range("D85:E85").Value = 0
$ws.Range("D86:E86").Value = 0
$ws.Range("D87:E87").Value = 0
$ws.Range("D88:E88").Value = 0
$ws.Range("D92:E92").Value = 0
$ws.Range("D93:E93").Value = 0
$ws.Range("D96:E96").Value = 0
$ws.Range("D97:E97").Value = 0
$ws.Range("D98:E98").Value = 0
$ws.Range("D99:E99").Value = 0
$ws.Range("D101:E101").Value = 0
